# Swap the identifying values between row 2 and row 4 on the "Artfynd" sheet.
# Columns affected: A (Id), Q (Ost), R (Nord), Z (Starttid), AB (Sluttid)
# Note: Value2 is used instead of Value because this runtime's Value
# property accessor does not reliably marshal scalar results back through
# PowerShell (it can yield a reflection artifact rather than the real data).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# --- Column A (Id) ---
$tmpA = $ws.Range("A2").Value2
$ws.Range("A2").Value2 = $ws.Range("A4").Value2
$ws.Range("A4").Value2 = $tmpA

# --- Column Q (Ost) ---
$tmpQ = $ws.Range("Q2").Value2
$ws.Range("Q2").Value2 = $ws.Range("Q4").Value2
$ws.Range("Q4").Value2 = $tmpQ

# --- Column R (Nord) ---
$tmpR = $ws.Range("R2").Value2
$ws.Range("R2").Value2 = $ws.Range("R4").Value2
$ws.Range("R4").Value2 = $tmpR

# --- Column Z (Starttid) ---
$tmpZ = $ws.Range("Z2").Value2
$ws.Range("Z2").Value2 = $ws.Range("Z4").Value2
$ws.Range("Z4").Value2 = $tmpZ

# --- Column AB (Sluttid) ---
$tmpAB = $ws.Range("AB2").Value2
$ws.Range("AB2").Value2 = $ws.Range("AB4").Value2
$ws.Range("AB4").Value2 = $tmpAB
